$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for rows 2-13 (columns D, K, L, M, N, O, P, Q, R, S, T)
# This reflects the weekly re-ordering of the Cereza price records.
$data = @{
    2  = @{ D = 44532; K = "Brooks";      L = "Primera"; M = 400; N = 27000; O = 28000; P = 27500; Q = "`$/bandeja 12 kilos"; R = "Región de O'Higgins"; S = 2292; T = 12 }
    3  = @{ D = 44557; K = "Lapins";      L = "Primera"; M = 250; N = 9000;  O = 10000; P = 9500;  Q = "`$/bandeja 10 kilos"; R = "Provincia de Curicó";  S = 950;  T = 10 }
    4  = @{ D = 44229; K = "Santina";     L = "Primera"; M = 250; N = 6500;  O = 7000;  P = 6750;  Q = "`$/bandeja 5 kilos";  R = "Provincia de Curicó";  S = 1350; T = 5  }
    5  = @{ D = 44210; K = "Rainier";     L = "Segunda"; M = 250; N = 21000; O = 22000; P = 21500; Q = "`$/caja 18 kilos";    R = "Región de O'Higgins"; S = 1194; T = 18 }
    6  = @{ D = 44175; K = "Rainier";     L = "Segunda"; M = 270; N = 25000; O = 26000; P = 25500; Q = "`$/caja 18 kilos";    R = "Región de O'Higgins"; S = 1417; T = 18 }
    7  = @{ D = 44208; K = "Lapins";      L = "Segunda"; M = 200; N = 10500; O = 11000; P = 10750; Q = "`$/bandeja 12 kilos"; R = "Provincia de Curicó";  S = 896;  T = 12 }
    8  = @{ D = 44568; K = "Santina";     L = "Segunda"; M = 200; N = 15000; O = 16000; P = 15500; Q = "`$/bandeja 12 kilos"; R = "Región de O'Higgins"; S = 1292; T = 12 }
    9  = @{ D = 44571; K = "Brooks";      L = "Segunda"; M = 400; N = 8500;  O = 9000;  P = 8750;  Q = "`$/bandeja 10 kilos"; R = "Región de O'Higgins"; S = 875;  T = 10 }
    10 = @{ D = 44161; K = "Bing";        L = "Primera"; M = 160; N = 39000; O = 40000; P = 39500; Q = "`$/caja 20 kilos";    R = "Provincia de Curicó";  S = 1975; T = 20 }
    11 = @{ D = 44580; K = "Sweet Heart"; L = "Segunda"; M = 300; N = 7000;  O = 8000;  P = 7500;  Q = "`$/bandeja 10 kilos"; R = "Región de O'Higgins"; S = 750;  T = 10 }
    12 = @{ D = 44594; K = "Santina";     L = "Primera"; M = 160; N = 5000;  O = 6000;  P = 5500;  Q = "`$/bandeja 5 kilos";  R = "Región de O'Higgins"; S = 1100; T = 5  }
    13 = @{ D = 44537; K = "Brooks";      L = "Primera"; M = 200; N = 29000; O = 30000; P = 29500; Q = "`$/caja 20 kilos";    R = "Región de O'Higgins"; S = 1475; T = 20 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("K$row").Value = $vals.K
    $ws.Range("L$row").Value = $vals.L
    $ws.Range("M$row").Value = $vals.M
    $ws.Range("N$row").Value = $vals.N
    $ws.Range("O$row").Value = $vals.O
    $ws.Range("P$row").Value = $vals.P
    $ws.Range("Q$row").Value = $vals.Q
    $ws.Range("R$row").Value = $vals.R
    $ws.Range("S$row").Value = $vals.S
    $ws.Range("T$row").Value = $vals.T
}
